# Generate Report for Archive
# Status of the localized items moved on from "Ready for handoff" to
# "In Translation" -- update the Overview rollup (zh-cn/de-de columns)
# as well as each language sheet's own Status column, then tighten the
# Status columns back up to the (now shorter) content width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: the zh-cn (col E) and de-de (col F) status cells.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# Per-language sheets: the Status column (col C).
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Re-fit the Status columns now that the text is shorter.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
